# Mapeo de objetos y ajuste para question de mensajes generica
#
# Adds a new "mensajeRespuesta" column (K) to the "Datos" sheet, fills in
# sample data for the existing numeroDocumento/usuario/clave/segundaClave
# columns on row 2, and gives the new header cell the same
# border/fill/font "pill" look used by the other header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- Row 2: sample data for numeroDocumento / usuario / clave / segundaClave ---
$ws.Range("F2").Value = "123456789"
$ws.Range("G2").Value = "prueba"
$ws.Range("H2").Value = 1234
$ws.Range("I2").Value = 1234

# --- New column K: "mensajeRespuesta" header + sample response text ---
$hdr = $ws.Range("K1")
$hdr.Value = "mensajeRespuesta"
$hdr.NumberFormat = "@"
$hdr.Font.Color = 0
$hdr.Font.Name = "Mic Shell Dlg"
$hdr.Interior.Pattern = 1
$hdr.Interior.Color = 5296274
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.Borders.Item(7).LineStyle = 1
$hdr.Borders.Item(10).LineStyle = 1

$ws.Range("K2").Value = "El usuario ha sido enviado al correo electronico"
$ws.Range("A2").Copy()
$ws.Range("K2").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

$ws.Columns.Item(11).ColumnWidth = 42.5

$ws.Application.CutCopyMode = $false

# --- Selection moves to the new header cell, matching the saved view ---
$ws.Range("K1").Select()
